# Updates cryptos list price/volume figures (and restores WEMIXToken/Filecoin
# row order) to match the latest GitHub Actions scrape.
#
# Note: several Price values (column D) are plain decimal numbers (e.g.
# "253.61", "0.700", "6.60"). Assigning such literal strings straight to
# .Value would make Excel auto-convert them to native numbers (and drop
# meaningful trailing zeros). To keep them as literal text - matching the
# original inlineStr/text cells - we temporarily force a text NumberFormat
# before writing the value, then restore the cell to the "Normal" style so
# no visible formatting/style changes are left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.160.74'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.905.09'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '253.61'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.78%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.700'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.85'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.57%  '
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '52.30'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('E11').Value = '  +5.11%  '
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.34'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.70%  '
$ws.Range('D14').Value = '2.181.22'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.736'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.00%  '
$ws.Range('E16').Value = '  +4.52%  '
$ws.Range('D17').Value = '1.904.19'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '35.165.89'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.82'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.46%  '
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '243.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '13.09'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('E23').Value = '  +5.35%  '
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.43'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.86%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.32'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '168.41'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = '4.128.74'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.35'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.48%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.04'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +9.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0598'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.36%  '
$ws.Range('E35').Value = '  +8.78%  '
$ws.Range('E36').Value = '  +3.42%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.851'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.34%  '
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.33'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.57%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '98.02'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.53%  '
$ws.Range('E42').Value = '  +4.08%  '
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0660'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '1.306.98'
$ws.Range('E46').Value = '  -3.01%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '12.41'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.60'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('E51').Value = '  +7.16%  '